$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 17.470401
$ws.Cells.Item(2, 8).Value = 52.411203
$ws.Cells.Item(2, 9).Value = 0.8600988665959021
$ws.Cells.Item(2, 10).Value = 0.8884442399952684
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 147.0592853333333
$ws.Cells.Item(2, 14).Value = 441.177856
$ws.Cells.Item(2, 15).Value = 0.9129893958419346
$ws.Cells.Item(2, 16).Value = 0.9274576550077637
$ws.Cells.Item(2, 17).Value = 2569.184685546752
$ws.Cells.Item(2, 18).Value = 23122.66216992077
$ws.Cells.Item(2, 19).Value = 0.7852611445777253
$ws.Cells.Item(2, 20).Value = 0.8239944114311665

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 17.470401
$ws.Cells.Item(3, 8).Value = 52.411203
$ws.Cells.Item(3, 9).Value = 0.8600988665959021
$ws.Cells.Item(3, 10).Value = 0.8884442399952684
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.3688046666666667
$ws.Cells.Item(3, 14).Value = 1.106414
$ws.Cells.Item(3, 15).Value = 0.00228965310854373
$ws.Cells.Item(3, 16).Value = 0.002325937532793486
$ws.Cells.Item(3, 17).Value = 6.443165417337999
$ws.Cells.Item(3, 18).Value = 57.988488756042
$ws.Cells.Item(3, 19).Value = 0.001969328043556246
$ws.Cells.Item(3, 20).Value = 0.002066465803599178

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 17.470401
$ws.Cells.Item(4, 8).Value = 52.411203
$ws.Cells.Item(4, 9).Value = 0.8600988665959021
$ws.Cells.Item(4, 10).Value = 0.8884442399952684
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 4.963579
$ws.Cells.Item(4, 14).Value = 14.890737
$ws.Cells.Item(4, 15).Value = 0.03081542917981618
$ws.Cells.Item(4, 16).Value = 0.03130376520837289
$ws.Cells.Item(4, 17).Value = 86.71571552517899
$ws.Cells.Item(4, 18).Value = 780.441439726611
$ws.Cells.Item(4, 19).Value = 0.02650431571122619
$ws.Cells.Item(4, 20).Value = 0.02781164988954318

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 17.470401
$ws.Cells.Item(5, 8).Value = 52.411203
$ws.Cells.Item(5, 9).Value = 0.8600988665959021
$ws.Cells.Item(5, 10).Value = 0.8884442399952684
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.14456
$ws.Cells.Item(5, 14).Value = 3.43368
$ws.Cells.Item(5, 15).Value = 0.007105781457704291
$ws.Cells.Item(5, 16).Value = 0.007218387680924443
$ws.Cells.Item(5, 17).Value = 19.99592216856
$ws.Cells.Item(5, 18).Value = 179.96329951704
$ws.Cells.Item(5, 19).Value = 0.006111674578049638
$ws.Cells.Item(5, 20).Value = 0.006413134957170125

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 17.470401
$ws.Cells.Item(6, 8).Value = 52.411203
$ws.Cells.Item(6, 9).Value = 0.8600988665959021
$ws.Cells.Item(6, 10).Value = 0.8884442399952684
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 7.5382435
$ws.Cells.Item(6, 14).Value = 15.076487
$ws.Cells.Item(6, 15).Value = 0.04679974041200103
$ws.Cells.Item(6, 16).Value = 0.0316942545701456
$ws.Cells.Item(6, 17).Value = 131.6961367806435
$ws.Cells.Item(6, 18).Value = 790.176820683861
$ws.Cells.Item(6, 19).Value = 0.04025240368534452
$ws.Cells.Item(6, 20).Value = 0.02815857791378957

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 0.8975426666666667
$ws.Cells.Item(7, 8).Value = 2.692628
$ws.Cells.Item(7, 9).Value = 0.04418761940962108
$ws.Cells.Item(7, 10).Value = 0.04564386429080782
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 147.0592853333333
$ws.Cells.Item(7, 14).Value = 441.177856
$ws.Cells.Item(7, 15).Value = 0.9129893958419346
$ws.Cells.Item(7, 16).Value = 0.9274576550077637
$ws.Cells.Item(7, 17).Value = 131.9919831161743
$ws.Cells.Item(7, 18).Value = 1187.927848045568
$ws.Cells.Item(7, 19).Value = 0.0403428279484833
$ws.Cells.Item(7, 20).Value = 0.04233275134064522

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 0.8975426666666667
$ws.Cells.Item(8, 8).Value = 2.692628
$ws.Cells.Item(8, 9).Value = 0.04418761940962108
$ws.Cells.Item(8, 10).Value = 0.04564386429080782
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.3688046666666667
$ws.Cells.Item(8, 14).Value = 1.106414
$ws.Cells.Item(8, 15).Value = 0.00228965310854373
$ws.Cells.Item(8, 16).Value = 0.002325937532793486
$ws.Cells.Item(8, 17).Value = 0.3310179239991111
$ws.Cells.Item(8, 18).Value = 2.979161315992
$ws.Cells.Item(8, 19).Value = 0.0001011743201403862
$ws.Cells.Item(8, 20).Value = 0.0001061647770957223

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 0.8975426666666667
$ws.Cells.Item(9, 8).Value = 2.692628
$ws.Cells.Item(9, 9).Value = 0.04418761940962108
$ws.Cells.Item(9, 10).Value = 0.04564386429080782
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 4.963579
$ws.Cells.Item(9, 14).Value = 14.890737
$ws.Cells.Item(9, 15).Value = 0.03081542917981618
$ws.Cells.Item(9, 16).Value = 0.03130376520837289
$ws.Cells.Item(9, 17).Value = 4.455023931870667
$ws.Cells.Item(9, 18).Value = 40.095215386836
$ws.Cells.Item(9, 19).Value = 0.001361660456541849
$ws.Cells.Item(9, 20).Value = 0.001428824810962283

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.8975426666666667
$ws.Cells.Item(10, 8).Value = 2.692628
$ws.Cells.Item(10, 9).Value = 0.04418761940962108
$ws.Cells.Item(10, 10).Value = 0.04564386429080782
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 1.14456
$ws.Cells.Item(10, 14).Value = 3.43368
$ws.Cells.Item(10, 15).Value = 0.007105781457704291
$ws.Cells.Item(10, 16).Value = 0.007218387680924443
$ws.Cells.Item(10, 17).Value = 1.02729143456
$ws.Cells.Item(10, 18).Value = 9.24562291104
$ws.Cells.Item(10, 19).Value = 0.0003139875666609798
$ws.Cells.Item(10, 20).Value = 0.0003294751077065543

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.8975426666666667
$ws.Cells.Item(11, 8).Value = 2.692628
$ws.Cells.Item(11, 9).Value = 0.04418761940962108
$ws.Cells.Item(11, 10).Value = 0.04564386429080782
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 7.5382435
$ws.Cells.Item(11, 14).Value = 15.076487
$ws.Cells.Item(11, 15).Value = 0.04679974041200103
$ws.Cells.Item(11, 16).Value = 0.0316942545701456
$ws.Cells.Item(11, 17).Value = 6.765895172972667
$ws.Cells.Item(11, 18).Value = 40.595371037836
$ws.Cells.Item(11, 19).Value = 0.002067969117794565
$ws.Cells.Item(11, 20).Value = 0.001446648254398041

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 7).Value = 1.944141
$ws.Cells.Item(12, 8).Value = 3.888282
$ws.Cells.Item(12, 9).Value = 0.09571351399447693
$ws.Cells.Item(12, 10).Value = 0.06591189571392365
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 147.0592853333333
$ws.Cells.Item(12, 14).Value = 441.177856
$ws.Cells.Item(12, 15).Value = 0.9129893958419346
$ws.Cells.Item(12, 16).Value = 0.9274576550077637
$ws.Cells.Item(12, 17).Value = 285.903986047232
$ws.Cells.Item(12, 18).Value = 1715.423916283392
$ws.Cells.Item(12, 19).Value = 0.08738542331572605
$ws.Cells.Item(12, 20).Value = 0.0611304922359519

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 7).Value = 1.944141
$ws.Cells.Item(13, 8).Value = 3.888282
$ws.Cells.Item(13, 9).Value = 0.09571351399447693
$ws.Cells.Item(13, 10).Value = 0.06591189571392365
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 0.3688046666666667
$ws.Cells.Item(13, 14).Value = 1.106414
$ws.Cells.Item(13, 15).Value = 0.00228965310854373
$ws.Cells.Item(13, 16).Value = 0.002325937532793486
$ws.Cells.Item(13, 17).Value = 0.717008273458
$ws.Cells.Item(13, 18).Value = 4.302049640748001
$ws.Cells.Item(13, 19).Value = 0.0002191507448470979
$ws.Cells.Item(13, 20).Value = 0.0001533069520985851

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 7).Value = 1.944141
$ws.Cells.Item(14, 8).Value = 3.888282
$ws.Cells.Item(14, 9).Value = 0.09571351399447693
$ws.Cells.Item(14, 10).Value = 0.06591189571392365
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 4.963579
$ws.Cells.Item(14, 14).Value = 14.890737
$ws.Cells.Item(14, 15).Value = 0.03081542917981618
$ws.Cells.Item(14, 16).Value = 0.03130376520837289
$ws.Cells.Item(14, 17).Value = 9.649897440639
$ws.Cells.Item(14, 18).Value = 57.899384643834
$ws.Cells.Item(14, 19).Value = 0.002949453012048149
$ws.Cells.Item(14, 20).Value = 0.002063290507867425

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 7).Value = 1.944141
$ws.Cells.Item(15, 8).Value = 3.888282
$ws.Cells.Item(15, 9).Value = 0.09571351399447693
$ws.Cells.Item(15, 10).Value = 0.06591189571392365
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 1.14456
$ws.Cells.Item(15, 14).Value = 3.43368
$ws.Cells.Item(15, 15).Value = 0.007105781457704291
$ws.Cells.Item(15, 16).Value = 0.007218387680924443
$ws.Cells.Item(15, 17).Value = 2.22518602296
$ws.Cells.Item(15, 18).Value = 13.35111613776
$ws.Cells.Item(15, 19).Value = 0.0006801193129936744
$ws.Cells.Item(15, 20).Value = 0.0004757776160477631

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 7).Value = 1.944141
$ws.Cells.Item(16, 8).Value = 3.888282
$ws.Cells.Item(16, 9).Value = 0.09571351399447693
$ws.Cells.Item(16, 10).Value = 0.06591189571392365
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 7.5382435
$ws.Cells.Item(16, 14).Value = 15.076487
$ws.Cells.Item(16, 15).Value = 0.04679974041200103
$ws.Cells.Item(16, 16).Value = 0.0316942545701456
$ws.Cells.Item(16, 17).Value = 14.6554082563335
$ws.Cells.Item(16, 18).Value = 58.62163302533401
$ws.Cells.Item(16, 19).Value = 0.004479367608861948
$ws.Cells.Item(16, 20).Value = 0.002089028401957985
